$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.21875
$ws.Range("C2").Value = 0.5104166666666666
$ws.Range("J2").Value = 0.02083333333333333
$ws.Range("P2").Value = 0.1875
$ws.Range("S2").Value = 0.0625
# Row 3
$ws.Range("B3").Value = 0.0196078431372549
$ws.Range("C3").Value = 0.0392156862745098
$ws.Range("J3").Value = 0.02941176470588235
$ws.Range("P3").Value = 0.7156862745098039
$ws.Range("S3").Value = 0.196078431372549
# Row 6
$ws.Range("B6").Value = 0.02717391304347826
$ws.Range("D6").Value = 0.0108695652173913
$ws.Range("E6").Value = 0.0108695652173913
$ws.Range("F6").Value = 0.05978260869565218
$ws.Range("J6").Value = 0.2336956521739131
$ws.Range("O6").Value = 0.005434782608695652
$ws.Range("Q6").Value = 0.1630434782608696
$ws.Range("R6").Value = 0.1141304347826087
$ws.Range("S6").Value = 0.375
# Row 7
$ws.Range("B7").Value = 0.09473684210526316
$ws.Range("D7").Value = 0.02631578947368421
$ws.Range("F7").Value = 0.1
$ws.Range("J7").Value = 0.1578947368421053
$ws.Range("O7").Value = 0.005263157894736842
$ws.Range("Q7").Value = 0.1789473684210526
$ws.Range("R7").Value = 0.07894736842105263
$ws.Range("S7").Value = 0.3578947368421053
# Row 8
$ws.Range("B8").Value = 0.07799442896935933
$ws.Range("D8").Value = 0.01671309192200557
$ws.Range("F8").Value = 0.06963788300835655
$ws.Range("J8").Value = 0.1392757660167131
$ws.Range("O8").Value = 0.01392757660167131
$ws.Range("Q8").Value = 0.1532033426183844
$ws.Range("R8").Value = 0.1197771587743733
$ws.Range("S8").Value = 0.4094707520891365
# Row 9
$ws.Range("B9").Value = 0.06666666666666667
$ws.Range("D9").Value = 0.02083333333333333
$ws.Range("F9").Value = 0.07083333333333333
$ws.Range("J9").Value = 0.1125
$ws.Range("O9").Value = 0.03333333333333333
$ws.Range("Q9").Value = 0.1791666666666667
$ws.Range("R9").Value = 0.075
$ws.Range("S9").Value = 0.4416666666666667
# Row 10
$ws.Range("B10").Value = 0.0728110599078341
$ws.Range("D10").Value = 0.02211981566820276
$ws.Range("F10").Value = 0.06267281105990784
$ws.Range("J10").Value = 0.1216589861751152
$ws.Range("O10").Value = 0.009216589861751152
$ws.Range("Q10").Value = 0.1981566820276498
$ws.Range("R10").Value = 0.08018433179723503
$ws.Range("S10").Value = 0.4331797235023042
# Row 11
$ws.Range("G11").Value = 0.134020618556701
$ws.Range("J11").Value = 0.08934707903780069
$ws.Range("K11").Value = 0.1890034364261168
$ws.Range("L11").Value = 0.5601374570446735
$ws.Range("S11").Value = 0.0274914089347079
# Row 12
$ws.Range("F12").Value = 0.005917159763313609
$ws.Range("G12").Value = 0.7633136094674556
$ws.Range("J12").Value = 0.1834319526627219
$ws.Range("K12").Value = 0.005917159763313609
$ws.Range("L12").Value = 0.01775147928994083
$ws.Range("S12").Value = 0.02366863905325444
# Row 13
$ws.Range("G13").Value = 0.6923076923076923
$ws.Range("J13").Value = 0.3076923076923077
# Row 14
$ws.Range("G14").Value = 1
# Row 15
$ws.Range("F15").Value = 0.0160427807486631
$ws.Range("H15").Value = 0.1283422459893048
$ws.Range("I15").Value = 0.1176470588235294
$ws.Range("J15").Value = 0.4117647058823529
$ws.Range("K15").Value = 0.0427807486631016
$ws.Range("M15").Value = 0.0106951871657754
$ws.Range("O15").Value = 0.0481283422459893
$ws.Range("S15").Value = 0.2245989304812834
# Row 16
$ws.Range("F16").Value = 0.007633587786259542
$ws.Range("H16").Value = 0.1679389312977099
$ws.Range("I16").Value = 0.0916030534351145
$ws.Range("J16").Value = 0.3511450381679389
$ws.Range("K16").Value = 0.1526717557251908
$ws.Range("M16").Value = 0.01526717557251908
$ws.Range("O16").Value = 0.05343511450381679
$ws.Range("S16").Value = 0.1603053435114504
# Row 17
$ws.Range("F17").Value = 0.01329787234042553
$ws.Range("H17").Value = 0.2021276595744681
$ws.Range("I17").Value = 0.1409574468085106
$ws.Range("J17").Value = 0.3882978723404255
$ws.Range("K17").Value = 0.101063829787234
$ws.Range("M17").Value = 0.01329787234042553
$ws.Range("O17").Value = 0.05319148936170213
$ws.Range("S17").Value = 0.08776595744680851
# Row 18
$ws.Range("F18").Value = 0.01639344262295082
$ws.Range("H18").Value = 0.1639344262295082
$ws.Range("I18").Value = 0.09836065573770492
$ws.Range("J18").Value = 0.3879781420765027
$ws.Range("K18").Value = 0.1748633879781421
$ws.Range("M18").Value = 0.00546448087431694
$ws.Range("O18").Value = 0.0546448087431694
$ws.Range("S18").Value = 0.09836065573770492
# Row 19
$ws.Range("F19").Value = 0.01155555555555556
$ws.Range("H19").Value = 0.192
$ws.Range("I19").Value = 0.1226666666666667
$ws.Range("J19").Value = 0.3502222222222222
$ws.Range("K19").Value = 0.1217777777777778
$ws.Range("M19").Value = 0.02577777777777778
$ws.Range("N19").Value = 0.001777777777777778
$ws.Range("O19").Value = 0.07733333333333334
$ws.Range("S19").Value = 0.09688888888888889

Write-Output "Applied 106 cell updates"